$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C3").Value = -13.51079999999999
$ws.Range("E3").Value = 16.33199999999999
$ws.Range("B12").Value = 5.146499999999998
$ws.Range("C14").Value = -12.87069999999999
$ws.Range("E20").Value = 16.05639999999999
$ws.Range("E25").Value = 17.16750000000001
$ws.Range("C26").Value = -11.3485
$ws.Range("B27").Value = 6.433400000000005
$ws.Range("E30").Value = 15.4587
$ws.Range("C31").Value = -13.648
$ws.Range("B32").Value = 6.488399999999997
$ws.Range("C35").Value = -12.1566
$ws.Range("B36").Value = 9.300600000000003
$ws.Range("C37").Value = -13.3239
$ws.Range("B38").Value = 5.430400000000002
$ws.Range("E44").Value = 16.65290000000001
$ws.Range("C45").Value = -13.49319999999999
$ws.Range("B46").Value = 7.491500000000002
$ws.Range("E47").Value = 16.2368
$ws.Range("C52").Value = -10.7771
$ws.Range("B54").Value = 4.527699999999998
$ws.Range("B55").Value = 5.066499999999998
$ws.Range("B56").Value = 4.890199999999999
$ws.Range("C57").Value = -14.37039999999999
$ws.Range("E58").Value = 16.5747
$ws.Range("B67").Value = 5.979099999999995
$ws.Range("B69").Value = 5.806599999999992
$ws.Range("B72").Value = 5.402600000000004
$ws.Range("E78").Value = 16.64660000000002
$ws.Range("C81").Value = -12.93
$ws.Range("B83").Value = 4.978099999999995
$ws.Range("C83").Value = -14.27629999999999
$ws.Range("E84").Value = 16.70439999999999
$ws.Range("B86").Value = 4.999300000000003
$ws.Range("E89").Value = 17.27460000000002
$ws.Range("B91").Value = 5.573600000000002
$ws.Range("E91").Value = 17.94680000000002
$ws.Range("E92").Value = 18.02000000000003
$ws.Range("B93").Value = 7.068799999999997
$ws.Range("E96").Value = 15.8053
$ws.Range("B99").Value = 4.475799999999998
$ws.Range("C100").Value = -13.05899999999999
$ws.Range("C102").Value = -12.08150000000001
$ws.Range("E102").Value = 16.7456
